$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.182.85"
$ws.Range("E2").Value = "  -4.94%  "

$ws.Range("D3").Value = "2.237.27"
$ws.Range("E3").Value = "  -5.67%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "319.10"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").Value = "101.19"
$ws.Range("E6").Value = "  -6.71%  "

$ws.Range("D7").Value = "0.588"
$ws.Range("E7").Value = "  -7.01%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  -8.08%  "

$ws.Range("D10").Value = "37.08"
$ws.Range("E10").Value = "  -9.41%  "

$ws.Range("D11").Value = "54.46"
$ws.Range("E11").Value = "  -2.76%  "

$ws.Range("D13").Value = "7.70"
$ws.Range("E13").Value = "  -9.62%  "

$ws.Range("E14").Value = "  -1.60%  "

$ws.Range("D15").Value = "2.576.24"
$ws.Range("E15").Value = "  -5.71%  "

$ws.Range("D16").Value = "0.864"
$ws.Range("E16").Value = "  -12.26%  "

$ws.Range("D17").Value = "14.36"
$ws.Range("E17").Value = "  -6.58%  "

$ws.Range("D18").Value = "2.234.41"
$ws.Range("E18").Value = "  -5.69%  "

$ws.Range("D19").Value = "43.108.65"
$ws.Range("E19").Value = "  -5.04%  "

$ws.Range("D20").Value = "14.31"
$ws.Range("E20").Value = "  -8.22%  "

$ws.Range("D21").Value = "0.0₃0967"
$ws.Range("E21").Value = "  -8.92%  "

$ws.Range("D22").Value = "6.54"
$ws.Range("E22").Value = "  -10.37%  "

$ws.Range("D23").Value = "65.34"
$ws.Range("E23").Value = "  -11.00%  "

$ws.Range("D24").Value = "3.19"
$ws.Range("E24").Value = "  -11.30%  "

$ws.Range("D25").Value = "238.64"
$ws.Range("E25").Value = "  -8.74%  "

$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -8.07%  "

$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("E28").Value = "  +2.31%  "

$ws.Range("D29").Value = "10.06"
$ws.Range("E29").Value = "  -9.85%  "

$ws.Range("E30").Value = "  -2.50%  "

$ws.Range("D31").Value = "6.41"
$ws.Range("E31").Value = "  -14.12%  "

$ws.Range("D32").Value = "35.59"
$ws.Range("E32").Value = "  -4.27%  "

$ws.Range("D33").Value = "20.52"
$ws.Range("E33").Value = "  -8.14%  "

$ws.Range("D34").Value = "0.0877"
$ws.Range("E34").Value = "  -9.23%  "

$ws.Range("D35").Value = "152.88"
$ws.Range("E35").Value = "  -8.28%  "

$ws.Range("E36").Value = "  -4.99%  "

$ws.Range("D37").Value = "3.12"
$ws.Range("E37").Value = "  +5.96%  "

$ws.Range("D38").Value = "1.95"
$ws.Range("E38").Value = "  +2.78%  "

$ws.Range("E39").Value = "  -6.90%  "

$ws.Range("D40").Value = "4.46"
$ws.Range("E40").Value = "  -5.28%  "

$ws.Range("D41").Value = "0.105"
$ws.Range("E41").Value = "  -10.82%  "

$ws.Range("D42").Value = "3.71"
$ws.Range("E42").Value = "  -5.95%  "

$ws.Range("E43").Value = "  -8.29%  "

$ws.Range("D44").Value = "12.96"
$ws.Range("E44").Value = "  -0.66%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").Value = "1.809.01"
$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").Value = "87.64"
$ws.Range("E47").Value = "  -11.30%  "

$ws.Range("D48").Value = "0.207"
$ws.Range("E48").Value = "  -9.11%  "

$ws.Range("E49").Value = "  -7.58%  "

$ws.Range("E50").Value = "  -10.33%  "

$ws.Range("D51").Value = "59.21"
$ws.Range("E51").Value = "  -16.00%  "
